# adding code repo location to planning ppt
#
# 1. Slide 6 ("Technical Aspects") / "TextBox 2": add a new first bullet
#    paragraph with the code-repo URL, split across several runs (the
#    hyperlink-looking pieces end up as their own runs), and resize/
#    reposition the textbox to its new autofit-grown extent.
# 2. Slide 8 ("Goals") / "TextBox 5": the "Deploy the application on
#    multiple clusters using AWS" bullet used to be split across two
#    runs - merge it back into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 6: "Technical Aspects" -> TextBox 2 (shape #2)
# ---------------------------------------------------------------------
$s6   = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6  = $shp6.TextFrame.TextRange

# Push a new empty paragraph in front of the existing first paragraph
# ("In Scala: "), then fill it in run-by-run so the URL pieces land in
# their own <a:r> elements, matching how PowerPoint splits runs around
# autocorrected / flagged tokens.
$firstPara = $tr6.Paragraphs(1, 1)
[void]$firstPara.InsertBefore(" `r")

$newPara = $tr6.Paragraphs(1, 1)
$newPara.Text = "Code Repo: https://"
[void]$newPara.InsertAfter("github.com")
[void]$newPara.InsertAfter("/")
[void]$newPara.InsertAfter("saravanvadivel")
[void]$newPara.InsertAfter("/")
[void]$newPara.InsertAfter("ScalaFinalProject")
[void]$newPara.InsertAfter(" ")

# Resize/reposition the textbox *after* the text edit so the shape's
# spAutoFit recalculation from the inserted paragraph doesn't get
# clobbered by this explicit sizing.
$shp6.Left   = 50.87484415
$shp6.Top    = 59.131456323
$shp6.Width  = 731.453491567
$shp6.Height = 501.651641803

# ---------------------------------------------------------------------
# Slide 8: "Goals" -> TextBox 5 (shape #3)
# ---------------------------------------------------------------------
$s8   = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(3)
$tr8  = $shp8.TextFrame.TextRange

# Paragraph 3 is "Deploy the application on multiple clusters using AWS",
# currently split into two runs: "Deploy the application on multiple
# clusters " + "using AWS". Merge them into a single run by deleting the
# second run's text and appending it onto the first run.
$deployPara = $tr8.Paragraphs(3, 1)
$run1Len    = ("Deploy the application on multiple clusters ").Length

$tail     = $deployPara.Characters($run1Len + 1, $deployPara.Length - $run1Len)
$tailText = $tail.Text
[void]$tail.Delete()

$headRun = $deployPara.Characters(1, $run1Len)
$headRun.Text = $headRun.Text + $tailText
